$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.588.06'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.18%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.335.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.63%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.50%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.24%  '

# Row 7
$ws.Range("E7").Value = '  -0.26%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("E9").Value = '  -4.11%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.65'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.96%  '

# Row 11
$ws.Range("E11").Value = '  -4.35%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.912.48'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.60%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.136'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.98%  '

# Row 14
$ws.Range("E14").Value = '  -5.77%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '66.629.93'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.16%  '

# Row 16
$ws.Range("E16").Value = '  -2.61%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.339.34'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.22%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '436.89'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.00%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.85%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.57'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.87%  '

# Row 21
$ws.Range("E21").Value = '  -3.05%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.55'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.48%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.06%  '

# Row 24
$ws.Range("E24").Value = '  -0.78%  '

# Row 25
$ws.Range("E25").Value = '  -4.55%  '

# Row 26
$ws.Range("E26").Value = '  -0.02%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.13%  '

# Row 28
$ws.Range("E28").Value = '  -0.06%  '

# Row 29
$ws.Range("E29").Value = '  -1.56%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.83'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.67%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.29'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.51%  '

# Row 33
$ws.Range("E33").Value = '  -2.95%  '

# Row 34
$ws.Range("E34").Value = '  -4.69%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '161.86'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.97%  '

# Row 36
$ws.Range("E36").Value = '  -4.88%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '27.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.24%  '

# Row 38
$ws.Range("E38").Value = '  -5.93%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.814.80'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.22%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.791'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.56%  '

# Row 41
$ws.Range("E41").Value = '  -4.04%  '

# Row 42
$ws.Range("E42").Value = '  -5.70%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.16'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.72%  '

# Row 44
$ws.Range("E44").Value = '  -3.74%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.39%  '

# Row 46
$ws.Range("E46").Value = '  -7.37%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '321.28'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.95%  '

# Row 48
$ws.Range("E48").Value = '  -4.01%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.980'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.46%  '

# Row 50
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.101'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.95%  '

# Row 51
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.17'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.89%  '
